# Auto-generated Excel COM-interop script applying a scheduled market-data
# refresh to the Excalibur_Profits workbook. For every affected Leve row this
# updates the currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ
# (H/I/J), LevePriceNQ / LevePriceHQ (K/L) and the resulting LeveProfitNQ / HQ
# (M/N) columns, per commit "chore: update Sheets via scheduled runner".

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 962.5
$ws.Range("J12").Value = 1200
$ws.Range("L12").Value = 1200
$ws.Range("N12").Value = -1540
$ws.Range("H86").Value = 3684.6956
$ws.Range("I86").Value = 2064
$ws.Range("K86").Value = 2064
$ws.Range("M86").Value = -941
$ws.Range("H88").Value = 5189.8423
$ws.Range("I88").Value = 3466
$ws.Range("J88").Value = 6741.3
$ws.Range("K88").Value = 3466
$ws.Range("L88").Value = 6741.3
$ws.Range("M88").Value = -3060
$ws.Range("N88").Value = -7553.3
$ws.Range("H89").Value = 3684.6956
$ws.Range("I89").Value = 2064
$ws.Range("K89").Value = 10320
$ws.Range("M89").Value = -4704
$ws.Range("H91").Value = 5189.8423
$ws.Range("I91").Value = 3466
$ws.Range("J91").Value = 6741.3
$ws.Range("K91").Value = 3466
$ws.Range("L91").Value = 6741.3
$ws.Range("M91").Value = -2062
$ws.Range("N91").Value = -9549.299999999999
$ws.Range("H112").Value = 1200.8334
$ws.Range("J112").Value = 1214.3846
$ws.Range("L112").Value = 3643.1538
$ws.Range("N112").Value = -5859.1538
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()  # row 126: cell removed (was -89880)
$ws.Range("H131").Value = 3238.3572
$ws.Range("I131").Value = 2193.2
$ws.Range("K131").Value = 6579.599999999999
$ws.Range("M131").Value = -1539.599999999999
$ws.Range("H135").Value = 5068.3
$ws.Range("I135").Value = 1087.7142
$ws.Range("K135").Value = 9789.427799999999
$ws.Range("M135").Value = -7254.427799999999
$ws.Range("H138").Value = 2701.9167
$ws.Range("I138").Value = 1022.7857
$ws.Range("K138").Value = 3068.3571
$ws.Range("M138").Value = 2071.6429

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3350.875
$ws.Range("I2").Value = 2468
$ws.Range("J2").Value = 5999.5
$ws.Range("K2").Value = 2468
$ws.Range("L2").Value = 5999.5
$ws.Range("M2").Value = -2355
$ws.Range("N2").Value = -6225.5
$ws.Range("H25").Value = 395
$ws.Range("I25").Value = 395
$ws.Range("K25").Value = 395
$ws.Range("M25").Value = 7
$ws.Range("H116").Value = 3350.875
$ws.Range("I116").Value = 2468
$ws.Range("J116").Value = 5999.5
$ws.Range("K116").Value = 2468
$ws.Range("L116").Value = 5999.5
$ws.Range("M116").Value = -174
$ws.Range("N116").Value = -10587.5
$ws.Range("H132").Value = 3732.9285
$ws.Range("I132").Value = 3553.1365
$ws.Range("J132").Value = 4392.1665
$ws.Range("K132").Value = 10659.4095
$ws.Range("L132").Value = 13176.4995
$ws.Range("M132").Value = -8129.4095
$ws.Range("N132").Value = -18236.4995

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3350.875
$ws.Range("I3").Value = 2468
$ws.Range("J3").Value = 5999.5
$ws.Range("K3").Value = 2468
$ws.Range("L3").Value = 5999.5
$ws.Range("M3").Value = -2354
$ws.Range("N3").Value = -6227.5
$ws.Range("H8").Value = 2224.8333
$ws.Range("I8").Value = 1087.5
$ws.Range("J8").Value = 4499.5
$ws.Range("K8").Value = 1087.5
$ws.Range("L8").Value = 4499.5
$ws.Range("M8").Value = -947.5
$ws.Range("N8").Value = -4779.5
$ws.Range("H20").Value = 1549.3529
$ws.Range("I20").Value = 1074.5555
$ws.Range("K20").Value = 1074.5555
$ws.Range("M20").Value = -827.5554999999999
$ws.Range("H134").Value = 2387.182
$ws.Range("I134").Value = 1425.7646
$ws.Range("J134").Value = 5656
$ws.Range("K134").Value = 4277.293799999999
$ws.Range("L134").Value = 16968
$ws.Range("M134").Value = -1742.293799999999
$ws.Range("N134").Value = -22038

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 10923
$ws.Range("I41").Value = 3307.8
$ws.Range("J41").Value = 48999
$ws.Range("K41").Value = 3307.8
$ws.Range("L41").Value = 48999
$ws.Range("M41").Value = -2879.8
$ws.Range("N41").Value = -49855
$ws.Range("H58").Value = 2354.8333
$ws.Range("I58").Value = 1370.7142
$ws.Range("J58").Value = 5799.25
$ws.Range("K58").Value = 1370.7142
$ws.Range("L58").Value = 5799.25
$ws.Range("M58").Value = -1167.7142
$ws.Range("N58").Value = -6205.25
$ws.Range("H103").Value = 14690.833
$ws.Range("I103").Value = 14690.833
$ws.Range("K103").Value = 14690.833
$ws.Range("M103").Value = -13518.833
$ws.Range("H122").Value = 2059.5557
$ws.Range("I122").Value = 1567
$ws.Range("K122").Value = 4701
$ws.Range("M122").Value = -2251
$ws.Range("H132").Value = 111116210
$ws.Range("J132").Value = 7014
$ws.Range("L132").Value = 21042
$ws.Range("N132").Value = -26102
$ws.Range("H134").Value = 13794.538
$ws.Range("J134").Value = 1822.1666
$ws.Range("L134").Value = 5466.4998
$ws.Range("N134").Value = -10536.4998
$ws.Range("H136").Value = 2354.8333
$ws.Range("I136").Value = 1370.7142
$ws.Range("J136").Value = 5799.25
$ws.Range("K136").Value = 4112.142599999999
$ws.Range("L136").Value = 17397.75
$ws.Range("M136").Value = -1562.142599999999
$ws.Range("N136").Value = -22497.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 647.5
$ws.Range("I5").Value = 647.5
$ws.Range("K5").Value = 1942.5
$ws.Range("M5").Value = -1830.5
$ws.Range("H133").Value = 8125.7144
$ws.Range("I133").Value = 5888.5
$ws.Range("J133").Value = 8498.583000000001
$ws.Range("K133").Value = 17665.5
$ws.Range("L133").Value = 25495.749
$ws.Range("M133").Value = -12605.5
$ws.Range("N133").Value = -35615.749
$ws.Range("H135").Value = 647.5
$ws.Range("I135").Value = 647.5
$ws.Range("K135").Value = 5827.5
$ws.Range("M135").Value = -3292.5
$ws.Range("H137").Value = 3174.6191
$ws.Range("J137").Value = 4036.1538
$ws.Range("L137").Value = 12108.4614
$ws.Range("N137").Value = -22308.4614
$ws.Range("H141").Value = 2629.5
$ws.Range("I141").Value = 2629.5
$ws.Range("K141").Value = 7888.5
$ws.Range("M141").Value = -2708.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 6906.1
$ws.Range("I9").Value = 880.2857
$ws.Range("J9").Value = 20966.334
$ws.Range("K9").Value = 880.2857
$ws.Range("L9").Value = 20966.334
$ws.Range("M9").Value = -710.2857
$ws.Range("N9").Value = -21306.334
$ws.Range("H80").Value = 1202120
$ws.Range("I80").Value = 5000000
$ws.Range("K80").Value = 5000000
$ws.Range("M80").Value = -4999002
$ws.Range("H83").Value = 1202120
$ws.Range("I83").Value = 5000000
$ws.Range("K83").Value = 25000000
$ws.Range("M83").Value = -24995008
$ws.Range("H97").Value = 2208.7144
$ws.Range("I97").Value = 2155.2727
$ws.Range("K97").Value = 2155.2727
$ws.Range("M97").Value = -1659.2727
$ws.Range("H126").Value = 2145.25
$ws.Range("I126").Value = 2037.7059
$ws.Range("K126").Value = 6113.1177
$ws.Range("M126").Value = -3643.1177

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5091.846
$ws.Range("I7").Value = 2381.4546
$ws.Range("J7").Value = 19999
$ws.Range("K7").Value = 2381.4546
$ws.Range("L7").Value = 19999
$ws.Range("M7").Value = -2269.4546
$ws.Range("N7").Value = -20223
$ws.Range("H61").Value = 1268.7059
$ws.Range("I61").Value = 871.7273
$ws.Range("J61").Value = 1996.5
$ws.Range("K61").Value = 871.7273
$ws.Range("L61").Value = 1996.5
$ws.Range("M61").Value = -669.7273
$ws.Range("N61").Value = -2400.5
$ws.Range("H82").Value = 2510.5715
$ws.Range("I82").Value = 1896.8334
$ws.Range("J82").Value = 2970.875
$ws.Range("K82").Value = 1896.8334
$ws.Range("L82").Value = 2970.875
$ws.Range("M82").Value = -1535.8334
$ws.Range("N82").Value = -3692.875
$ws.Range("H85").Value = 2510.5715
$ws.Range("I85").Value = 1896.8334
$ws.Range("J85").Value = 2970.875
$ws.Range("K85").Value = 1896.8334
$ws.Range("L85").Value = 2970.875
$ws.Range("M85").Value = -648.8334
$ws.Range("N85").Value = -5466.875
$ws.Range("H100").Value = 13724.728
$ws.Range("I100").Value = 4774.222
$ws.Range("K100").Value = 4774.222
$ws.Range("M100").Value = -4233.222
$ws.Range("H113").Value = 1268.7059
$ws.Range("I113").Value = 871.7273
$ws.Range("J113").Value = 1996.5
$ws.Range("K113").Value = 871.7273
$ws.Range("L113").Value = 1996.5
$ws.Range("M113").Value = 1298.2727
$ws.Range("N113").Value = -6336.5
$ws.Range("H126").Value = 5091.846
$ws.Range("I126").Value = 2381.4546
$ws.Range("J126").Value = 19999
$ws.Range("K126").Value = 7144.3638
$ws.Range("L126").Value = 59997
$ws.Range("M126").Value = -4674.3638
$ws.Range("N126").Value = -64937
$ws.Range("H132").Value = 2960.4736
$ws.Range("J132").Value = 4332.3335
$ws.Range("L132").Value = 12997.0005
$ws.Range("N132").Value = -18057.0005

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5212.3335
$ws.Range("I81").Value = 5212.3335
$ws.Range("K81").Value = 10424.667
$ws.Range("M81").Value = -9363.666999999999
$ws.Range("H84").Value = 5212.3335
$ws.Range("I84").Value = 5212.3335
$ws.Range("K84").Value = 52123.335
$ws.Range("M84").Value = -46819.335
$ws.Range("H100").Value = 4876
$ws.Range("I100").Value = 3946.6
$ws.Range("K100").Value = 7893.2
$ws.Range("M100").Value = -7352.2
$ws.Range("H132").Value = 3393035.8
$ws.Range("I132").Value = 4168833.2
$ws.Range("J132").Value = 7737
$ws.Range("K132").Value = 12506499.6
$ws.Range("L132").Value = 23211
$ws.Range("M132").Value = -12503969.6
$ws.Range("N132").Value = -28271
$ws.Range("H136").Value = 5749417.5
$ws.Range("I136").Value = 6062786
$ws.Range("K136").Value = 18188358
$ws.Range("M136").Value = -18185808

